$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$excel.ActiveWindow.ScrollRow = 30
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A31").Select()
